$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2058047493403694
$ws.Range("C2").Value = 0.496042216358839
$ws.Range("J2").Value = 0.0237467018469657
$ws.Range("P2").Value = 0.1556728232189974
$ws.Range("S2").Value = 0.1187335092348285
$ws.Range("B3").Value = 0.02116402116402116
$ws.Range("C3").Value = 0.01587301587301587
$ws.Range("J3").Value = 0.07407407407407407
$ws.Range("P3").Value = 0.7248677248677249
$ws.Range("S3").Value = 0.164021164021164
$ws.Range("J4").Value = 0.05882352941176471
$ws.Range("P4").Value = 0.8529411764705882
$ws.Range("S4").Value = 0.08823529411764706
$ws.Range("B6").Value = 0.0735930735930736
$ws.Range("D6").Value = 0.004329004329004329
$ws.Range("F6").Value = 0.04329004329004329
$ws.Range("J6").Value = 0.2770562770562771
$ws.Range("O6").Value = 0.008658008658008658
$ws.Range("Q6").Value = 0.1991341991341991
$ws.Range("R6").Value = 0.06493506493506493
$ws.Range("S6").Value = 0.329004329004329
$ws.Range("B7").Value = 0.1010830324909747
$ws.Range("D7").Value = 0.01083032490974729
$ws.Range("F7").Value = 0.04693140794223827
$ws.Range("J7").Value = 0.1696750902527076
$ws.Range("O7").Value = 0.01805054151624549
$ws.Range("Q7").Value = 0.1624548736462094
$ws.Range("R7").Value = 0.08303249097472924
$ws.Range("S7").Value = 0.407942238267148
$ws.Range("B8").Value = 0.1138014527845036
$ws.Range("D8").Value = 0.01694915254237288
$ws.Range("F8").Value = 0.06053268765133172
$ws.Range("J8").Value = 0.1670702179176755
$ws.Range("O8").Value = 0.02663438256658596
$ws.Range("Q8").Value = 0.1307506053268765
$ws.Range("R8").Value = 0.05084745762711865
$ws.Range("S8").Value = 0.4334140435835351
$ws.Range("B9").Value = 0.1275510204081633
$ws.Range("D9").Value = 0.00510204081632653
$ws.Range("F9").Value = 0.05102040816326531
$ws.Range("J9").Value = 0.1479591836734694
$ws.Range("O9").Value = 0.01020408163265306
$ws.Range("Q9").Value = 0.1581632653061225
$ws.Range("R9").Value = 0.07653061224489796
$ws.Range("S9").Value = 0.4234693877551021
$ws.Range("B10").Value = 0.1151178918169209
$ws.Range("D10").Value = 0.01803051317614424
$ws.Range("E10").Value = 0.0006934812760055479
$ws.Range("F10").Value = 0.06657420249653259
$ws.Range("J10").Value = 0.1484049930651872
$ws.Range("O10").Value = 0.03467406380027739
$ws.Range("Q10").Value = 0.20249653259362
$ws.Range("R10").Value = 0.05825242718446602
$ws.Range("S10").Value = 0.3557558945908461
$ws.Range("G11").Value = 0.1543942992874109
$ws.Range("J11").Value = 0.09501187648456057
$ws.Range("K11").Value = 0.2019002375296912
$ws.Range("L11").Value = 0.5415676959619953
$ws.Range("S11").Value = 0.007125890736342043
$ws.Range("G12").Value = 0.7468879668049793
$ws.Range("J12").Value = 0.1784232365145228
$ws.Range("K12").Value = 0.01244813278008299
$ws.Range("L12").Value = 0.03734439834024896
$ws.Range("S12").Value = 0.02489626556016597
$ws.Range("F13").Value = 0.02083333333333333
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.2083333333333333
$ws.Range("S13").Value = 0.02083333333333333
$ws.Range("F15").Value = 0.02409638554216868
$ws.Range("H15").Value = 0.1726907630522088
$ws.Range("I15").Value = 0.07228915662650602
$ws.Range("J15").Value = 0.3493975903614458
$ws.Range("K15").Value = 0.08433734939759036
$ws.Range("M15").Value = 0.01606425702811245
$ws.Range("O15").Value = 0.05622489959839357
$ws.Range("S15").Value = 0.2248995983935743
$ws.Range("F16").Value = 0.009174311926605505
$ws.Range("H16").Value = 0.1330275229357798
$ws.Range("I16").Value = 0.06422018348623854
$ws.Range("J16").Value = 0.481651376146789
$ws.Range("K16").Value = 0.1009174311926606
$ws.Range("M16").Value = 0.03211009174311927
$ws.Range("O16").Value = 0.04128440366972477
$ws.Range("S16").Value = 0.1376146788990826
$ws.Range("F17").Value = 0.01108647450110865
$ws.Range("H17").Value = 0.1485587583148559
$ws.Range("I17").Value = 0.09977827050997783
$ws.Range("J17").Value = 0.4235033259423504
$ws.Range("K17").Value = 0.1352549889135255
$ws.Range("M17").Value = 0.02439024390243903
$ws.Range("O17").Value = 0.03991130820399113
$ws.Range("S17").Value = 0.1175166297117517
$ws.Range("F18").Value = 0.01282051282051282
$ws.Range("H18").Value = 0.1346153846153846
$ws.Range("I18").Value = 0.1217948717948718
$ws.Range("J18").Value = 0.4230769230769231
$ws.Range("K18").Value = 0.1602564102564103
$ws.Range("M18").Value = 0.01282051282051282
$ws.Range("O18").Value = 0.03205128205128205
$ws.Range("S18").Value = 0.1025641025641026
$ws.Range("F19").Value = 0.02686567164179104
$ws.Range("H19").Value = 0.1962686567164179
$ws.Range("I19").Value = 0.07388059701492537
$ws.Range("J19").Value = 0.3529850746268657
$ws.Range("K19").Value = 0.1432835820895522
$ws.Range("M19").Value = 0.01716417910447761
$ws.Range("N19").Value = 0.0007462686567164179
$ws.Range("O19").Value = 0.07985074626865672
$ws.Range("S19").Value = 0.108955223880597
